# Applies the changes described in the target diff:
#  1. Rename sheet 1 "Export this as TSV" -> "Export as TSV"
#  2. Rename sheet 8 / 9 (their names were re-truncated to 31 chars by the
#     generator that produced the target workbook)
#  3. Freeze the header row on the main worksheet
#  4. Add errorTitle/error text to every data validation rule on the main
#     worksheet, and repoint the two list validations whose source sheets
#     were renamed

$wb = $excel.ActiveWorkbook

# --- 1 & 2: sheet renames -------------------------------------------------
$wb.Worksheets.Item(1).Name = "Export as TSV"
$wb.Worksheets.Item(8).Name = "specimen_prese...mperature list"
$wb.Worksheets.Item(9).Name = "specimen_tumor...ance_unit list"

$ws = $wb.Worksheets.Item(1)

# --- 3: freeze header row -------------------------------------------------
$ws.Activate()
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)

# --- 4: data validations ---------------------------------------------------
$rng = $ws.Range("B2:B1048576")
$rng.Validation.Modify(3, 1, 1, "='vital_state list'!`$A`$1:`$A`$2")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: living / deceased."

$rng = $ws.Range("C2:C1048576")
$rng.Validation.Modify(3, 1, 1, "='health_status list'!`$A`$1:`$A`$3")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: cancer / relatively healthy / chronic illness."

$rng = $ws.Range("D2:D1048576")
$rng.Validation.Modify(3, 1, 1, "='organ_condition list'!`$A`$1:`$A`$2")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: healthy / diseased."

$rng = $ws.Range("F2:F1048576")
$rng.Validation.Modify(3, 1, 1, "='perfusion_solution list'!`$A`$1:`$A`$4")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: UWS / HTK / Unknown / None."

$rng = $ws.Range("H2:H1048576")
$rng.Validation.Modify(2, 1, 1, "-1e+307", "1e+307")
$rng.Validation.ErrorTitle = "Not a number"
$rng.Validation.ErrorMessage = "The values in this column must be numbers."

$rng = $ws.Range("I2:I1048576")
$rng.Validation.Modify(3, 1, 1, "='warm_ischemia_time_unit list'!`$A`$1:`$A`$1")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: minutes."

$rng = $ws.Range("J2:J1048576")
$rng.Validation.Modify(2, 1, 1, "-1e+307", "1e+307")
$rng.Validation.ErrorTitle = "Not a number"
$rng.Validation.ErrorMessage = "The values in this column must be numbers."

$rng = $ws.Range("K2:K1048576")
$rng.Validation.Modify(3, 1, 1, "='cold_ischemia_time_unit list'!`$A`$1:`$A`$1")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: minutes."

$rng = $ws.Range("L2:L1048576")
$rng.Validation.Modify(3, 1, 1, "='specimen_prese...mperature list'!`$A`$1:`$A`$5")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: Liquid Nitrogen / Liquid Nitrogen Vapor / Freezer (-80 Celsius) / Freezer (-20 Celsius) / Room Temperature."

$rng = $ws.Range("N2:N1048576")
$rng.Validation.Modify(2, 1, 1, "-1e+307", "1e+307")
$rng.Validation.ErrorTitle = "Not a number"
$rng.Validation.ErrorMessage = "The values in this column must be numbers."

$rng = $ws.Range("O2:O1048576")
$rng.Validation.Modify(3, 1, 1, "='specimen_tumor...ance_unit list'!`$A`$1:`$A`$1")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: cm."
